# Update countries & provincias Spain
# Applies the COVID data refresh captured in the target diff:
#  - bumps the "Datos actualizados" timestamp
#  - refreshes Kazajistan's row (row 32)
#  - Honduras overtakes Suiza in the ranking (rows 54 & 55 swap identity + get new figures)
#  - Islas Malvinas / Groenlandia swap places (tie on "Casos totales", same figures either way)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 20 de Julio de 2020 a las 05:48"

# 2) Kazajistan (row 32) - refreshed totals
$ws.Range("B32").Value = 71838
$ws.Range("C32").Value = 1499
$ws.Range("E32").Value = 28434

# 3) Honduras moves ahead of Suiza (row 54 becomes Honduras, row 55 becomes Suiza)
$ws.Range("A54").Value = "Honduras"
$ws.Range("B54").Value = 33835
$ws.Range("C54").Value = 1042
$ws.Range("D54").Value = 3801
$ws.Range("E54").Value = 29134
$ws.Range("F54").Value = 0
$ws.Range("G54").Value = 9
$ws.Range("H54").Value = 900

$ws.Range("A55").Value = "Suiza"
$ws.Range("B55").Value = 33591
$ws.Range("C55").Value = 0
$ws.Range("D55").Value = 30300
$ws.Range("E55").Value = 1322
$ws.Range("F55").Value = 0
$ws.Range("G55").Value = 0
$ws.Range("H55").Value = 1969

# 4) Islas Malvinas / Groenlandia swap order (figures are identical, only labels move)
$ws.Range("A210").Value = "Islas Malvinas"
$ws.Range("A211").Value = "Groenlandia"
